# Fill in the PowerPoint template: the subtitle on the title slide reads
# "Case Study 1 S S2020" and should read "Case Study 1 SE S2020".
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$subtitle = $s.Shapes.Item(2)
$subtitle.TextFrame.TextRange.Text = "Case Study 1 SE S2020"
